$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row 5 with test-case data (special characters test case)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "special characters"
$ws.Range("C5").Value = "# @ ass 1.2"
$ws.Range("D5").Value = "char 1 char 1 string 3 float 8"
$ws.Range("E5").Value = "char 1 char 1 string 3 float 8"
$ws.Range("F5").Value = "PASS"

# Update selection as in the diff
$ws.Range("J8").Select()
